# This workbook has 3 sheets: "БИВТ-22-17" (1), "БИВТ-22-18" (2), "БИВТ-22-20" (3).
# The edit touches the 3rd sheet ("БИВТ-22-20"), updating several lab-score cells,
# which in turn recalculates the dependent SUM / FLOOR.MATH / COUNT formulas, and
# moves the active tab/selection to that sheet at cell C25 (the last cell edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Row 2 (Акимов Дмитрий Дмитриевич): fill in the 5th lab score.
$ws.Range("E2").Value = 5

# Row 3 (Бондаревич Никита Валерьевич): fill in the 4th lab score.
$ws.Range("D3").Value = 5

# Row 7 (Гайдаржи Иван Владимирович): 4th lab score was "pass" -> now a numeric 5.
$ws.Range("D7").Value = 5

# Row 13 (Калашникова Анастасия Павловна): 5th lab score becomes "pass".
$ws.Range("E13").Value = "pass"

# Row 18 (Наумкин Павел Дмитриевич): 4th lab score was "pass" -> now a numeric 5.
$ws.Range("D18").Value = 5

# Row 25 (Тошматов Абдурахимжон Баходиржон угли): 2nd lab score "pass" -> "failed 1".
$ws.Range("C25").Value = "failed 1"

# Row 27 (Чавыкин Олег Сергеевич): 4th lab score was "pass" -> now a numeric 5.
$ws.Range("D27").Value = 5

# The last worked-on sheet/cell becomes the active tab/selection.
$ws.Activate() | Out-Null
$ws.Range("C25").Select() | Out-Null
